$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '22.381.18'
$ws.Range("E2").Value = '  -4.49%  '
$ws.Range("D3").Value = '1.570.34'
$ws.Range("E3").Value = '  -4.65%  '
$ws.Range("D4").Value = '1.001'
$ws.Range("E4").Value = '  -0.08%  '
$ws.Range("E5").Value = '  -0.13%  '
$ws.Range("D6").Value = '291.08'
$ws.Range("E6").Value = '  -2.79%  '
$ws.Range("D7").Value = '0.3663'
$ws.Range("E7").Value = '  -3.34%  '
$ws.Range("D8").Value = '49.26'
$ws.Range("E8").Value = '  -1.02%  '
$ws.Range("D9").Value = '0.3371'
$ws.Range("E9").Value = '  -5.21%  '
$ws.Range("D10").Value = '1.171'
$ws.Range("E10").Value = '  -4.04%  '
$ws.Range("D11").Value = '0.07581'
$ws.Range("E11").Value = '  -6.35%  '
$ws.Range("D12").Value = '1.001'
$ws.Range("E12").Value = '  -0.01%  '
$ws.Range("D13").Value = '21.07'
$ws.Range("E13").Value = '  -4.31%  '
$ws.Range("D14").Value = '6.047'
$ws.Range("E14").Value = '  -5.41%  '
$ws.Range("D15").Value = '6.867'
$ws.Range("E15").Value = '  -6.70%  '
$ws.Range("D16").Value = '0.00001142'
$ws.Range("E16").Value = '  -4.50%  '
$ws.Range("D17").Value = '1.568.51'
$ws.Range("E17").Value = '  -4.47%  '
$ws.Range("D18").Value = '89.04'
$ws.Range("E18").Value = '  -8.58%  '
$ws.Range("D19").Value = '0.06707'
$ws.Range("E19").Value = '  -3.38%  '
$ws.Range("E20").Value = '  -0.06%  '
$ws.Range("D21").Value = '6.261'
$ws.Range("E21").Value = '  -7.51%  '
$ws.Range("D22").Value = '16.39'
$ws.Range("E22").Value = '  -5.27%  '
$ws.Range("D23").Value = '0.5234'
$ws.Range("E23").Value = '  -9.12%  '
$ws.Range("E24").Value = '  -3.36%  '
$ws.Range("D25").Value = '22.402.89'
$ws.Range("E25").Value = '  -4.51%  '
$ws.Range("D26").Value = '2.381'
$ws.Range("E26").Value = '  -4.75%  '
$ws.Range("D27").Value = '2.988'
$ws.Range("E27").Value = '  +2.41%  '
$ws.Range("D28").Value = '19.84'
$ws.Range("E28").Value = '  -5.07%  '
$ws.Range("D29").Value = '145.17'
$ws.Range("E29").Value = '  -5.02%  '
$ws.Range("D30").Value = '4.953'
$ws.Range("E30").Value = '  -5.06%  '
$ws.Range("D31").Value = '124.93'
$ws.Range("E31").Value = '  -5.95%  '
$ws.Range("D32").Value = '1.745.79'
$ws.Range("E32").Value = '  -4.47%  '
$ws.Range("D33").Value = '6.267'
$ws.Range("E33").Value = '  -9.54%  '
$ws.Range("D34").Value = '1.009'
$ws.Range("E34").Value = '  +0.04%  '
$ws.Range("D35").Value = '1.971'
$ws.Range("E35").Value = '  -5.97%  '
$ws.Range("D36").Value = '10.36'
$ws.Range("E36").Value = '  -11.67%  '
$ws.Range("D37").Value = '0.08429'
$ws.Range("E37").Value = '  -3.44%  '
$ws.Range("D38").Value = '0.02545'
$ws.Range("E38").Value = '  -6.49%  '
$ws.Range("D39").Value = '0.2299'
$ws.Range("E39").Value = '  -5.48%  '
$ws.Range("D40").Value = '5.522'
$ws.Range("E40").Value = '  -6.79%  '
$ws.Range("D41").Value = '0.06502'
$ws.Range("E41").Value = '  -3.92%  '
$ws.Range("D42").Value = '11.77'
$ws.Range("E42").Value = '  -9.69%  '
$ws.Range("E43").Value = '  -3.91%  '
$ws.Range("D44").Value = '0.6390'
$ws.Range("E44").Value = '  -7.08%  '
$ws.Range("D45").Value = '14.52'
$ws.Range("E45").Value = '  -6.75%  '
$ws.Range("D46").Value = '0.9997'
$ws.Range("E46").Value = '  -0.13%  '
$ws.Range("D47").Value = '0.6012'
$ws.Range("E47").Value = '  -5.94%  '
$ws.Range("E48").Value = '  -3.85%  '
$ws.Range("D49").Value = '2.121'
$ws.Range("E49").Value = '  -5.68%  '
$ws.Range("D50").Value = '1.202'
$ws.Range("E50").Value = '  +2.38%  '
$ws.Range("D51").Value = '121.27'
$ws.Range("E51").Value = '  -4.68%  '
